$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.082.57"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.913.06"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.7912"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.86%  "
$ws.Range("D6").Value = "'243.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.3173"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("D9").Value = "'26.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'0.06946"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.07997"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.915.52"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7505"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "'5.233"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'93.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "30.099.70"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'14.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "'247.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").Value = "'0.000007796"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'6.924"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("D24").Value = "'168.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "'9.321"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'0.1390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.70%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "'2.059"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'1.382"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "'1.525"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "'4.346"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'0.05747"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'1.263"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").Value = "'0.7376"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D37").Value = "'0.01922"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "'2.795"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'6.180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'72.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.907"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "'0.8319"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'7.606"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'101.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'9.881"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'990.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.10%  "
$ws.Range("D49").Value = "2.064.54"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'36.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  +2.22%  "
